$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Andrea Favero
$ws.Range("B2").Formula = "=2+5"
$ws.Range("C2").Value = 16
$ws.Range("D2").Formula = "=2+1+5+3"
$ws.Range("E2").Formula = "=10+8+10+7"
$ws.Range("F2").Formula = "=4+15+5"
$ws.Range("G2").Formula = "=10+2+5+15+10"

# Row 3 - Eleonora Thiella
$ws.Range("B3").Formula = "=3+5"
$ws.Range("C3").Formula = "=6+3"
$ws.Range("D3").Formula = "=18+3"
$ws.Range("E3").Formula = "=12+5+20"
$ws.Range("F3").Formula = "=2+2+15+5"
$ws.Range("G3").Formula = "=4+2+8+3+3+17"

# Row 4 - Federico Caldart
$ws.Range("B4").Formula = "=2+4"
$ws.Range("C4").Formula = "=4+1+3+5"
$ws.Range("D4").Formula = "=6+4+5"
$ws.Range("E4").Formula = "=16+20"
$ws.Range("F4").Formula = "=1+4+16"
$ws.Range("G4").Formula = "=18+3+10+13"

# Row 5 - Giovanni Cavallin
$ws.Range("B5").Formula = "=2+8"
$ws.Range("C5").Formula = "=6+2+5"
$ws.Range("D5").Formula = "=16+3"
$ws.Range("E5").Formula = "=16+4+15"
$ws.Range("F5").Formula = "=1+2+14+3"
$ws.Range("G5").Formula = "=4+2+7+5+16+4"

# Row 6 - Giovanni Dalla Riva
$ws.Range("B6").Formula = "=2+2+3"
$ws.Range("C6").Formula = "=6+3"
$ws.Range("D6").Formula = "=18+3"
$ws.Range("E6").Formula = "=6+5+20+4"
$ws.Range("F6").Formula = "=6+2+10+3"
$ws.Range("G6").Formula = "=4+9+2+16+11"

# Row 7 - Lorenzo Menegon
$ws.Range("B7").Formula = "=3+5"
$ws.Range("C7").Formula = "=4+2+3"
$ws.Range("D7").Formula = "=2+8"
$ws.Range("E7").Formula = "=8+8+15+4"
$ws.Range("F7").Formula = "=5+15"
$ws.Range("G7").Formula = "=22+3+10+18"

# Row 8 - Stefano Panozzo
$ws.Range("B8").Formula = "=22"
$ws.Range("C8").Formula = "=4+2"
$ws.Range("D8").Formula = "=4+4+2+2"
$ws.Range("E8").Formula = "=12+5+20"
$ws.Range("F8").Formula = "=1+15+4"
$ws.Range("G8").Formula = "=2+1+6+2+10+17"

# Row 9 TOTALE - G9 is re-entered as its own (non-shared) SUM formula, which
# Excel naturally breaks out of the previous C9:H9 shared-formula group.
$ws.Range("G9").Formula = "=SUM(G2:G8)"

# G9 picks up a top border (in addition to its existing right border) when the
# total row is finished off.
$ws.Range("G9").Borders.Item(8).LineStyle = 1
$ws.Range("G9").Borders.Item(8).Weight = 2

# Move the active selection off H9 (no longer the last-touched cell).
$ws.Range("K1").Select()
